$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = 0.8021519022689512
$ws.Range("J19").Value = 0.260316539010564
$ws.Range("K19").Value = 0.2282781566817243
$ws.Range("L19").Value = 2.418057025656014
